# "add lattice protein analysis" -- insert a new "Non-linear score"
# column (F) between the existing "Non-Linear fit success" column (E)
# and the "Saturation A"/"Saturation B" columns (old F/G, now shifted to
# G/H), populate a few new data points, and center-align the
# "Significant Order" (D) / "Non-Linear fit success" (E) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F. This shifts the old "Saturation A" (F) and
# "Saturation B" (G) columns right to G/H, carrying their widths/values
# along with them.
$ws.Columns.Item(6).Insert()

# New column F: header + the one new data point (dataset in row 13).
$ws.Range("F1").Value = "Non-linear score"
$ws.Range("F13").Value = 0.96

# New numeric entries in column D ("Significant Order").
$ws.Range("D2").Value = 4
$ws.Range("D10").Value = 4
$ws.Range("D13").Value = 5

# Center-align column D's header + its (now three) populated cells, and
# column E's header + all of its cells.
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D10").HorizontalAlignment = -4108
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("E1:E13").HorizontalAlignment = -4108

# Give the new column a sensible width (closest achievable value to the
# source width under this host's column-width pixel grid).
$ws.Columns.Item(6).ColumnWidth = 18.5

# Match the saved selection state.
$ws.Range("D6").Select() | Out-Null
